{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the worksheet\n// table with their updated values. Each old string is unique in the\n// document, so a targeted search/replace per pair is unambiguous. The\n// pairs are applied in document order (top-to-bottom, left-to-right),\n// which also guarantees that a later pair's \"old\" text is never a\n// string that an earlier pair just produced (verified against the diff).\nconst pairs = [\n  [\"60\u00f74=15, 0\", \"56\u00f74=14, 0\"],\n  [\"72\u00f79=8, 0\", \"36\u00f75=7, 1\"],\n  [\"64\u00f72=32, 0\", \"59\u00f76=9, 5\"],\n  [\"73\u00f76=12, 1\", \"24\u00f76=4, 0\"],\n  [\"27\u00f77=3, 6\", \"28\u00f73=9, 1\"],\n  [\"93\u00f76=15, 3\", \"40\u00f78=5, 0\"],\n  [\"51\u00f75=10, 1\", \"72\u00f74=18, 0\"],\n  [\"81\u00f75=16, 1\", \"30\u00f79=3, 3\"],\n  [\"72\u00f75=14, 2\", \"46\u00f78=5, 6\"],\n  [\"83\u00f74=20, 3\", \"10\u00f72=5, 0\"],\n  [\"58\u00f72=29, 0\", \"41\u00f75=8, 1\"],\n  [\"89\u00f72=44, 1\", \"39\u00f72=19, 1\"],\n  [\"47\u00f78=5, 7\", \"50\u00f74=12, 2\"],\n  [\"94\u00f72=47, 0\", \"49\u00f77=7, 0\"],\n  [\"26\u00f75=5, 1\", \"56\u00f77=8, 0\"],\n  [\"97\u00f72=48, 1\", \"83\u00f77=11, 6\"],\n  [\"29\u00f72=14, 1\", \"46\u00f76=7, 4\"],\n  [\"32\u00f78=4, 0\", \"51\u00f76=8, 3\"],\n  [\"48\u00f75=9, 3\", \"85\u00f76=14, 1\"],\n  [\"54\u00f72=27, 0\", \"81\u00f77=11, 4\"],\n  [\"53\u00f76=8, 5\", \"58\u00f78=7, 2\"],\n  [\"86\u00f72=43, 0\", \"34\u00f74=8, 2\"],\n  [\"38\u00f73=12, 2\", \"93\u00f76=15, 3\"],\n  [\"87\u00f78=10, 7\", \"29\u00f78=3, 5\"],\n  [\"36\u00f76=6, 0\", \"85\u00f74=21, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the worksheet\n# table with their updated values. Each old string is unique in the\n# document, so Find/Replace per pair is unambiguous. Pairs are applied in\n# document order (top-to-bottom, left-to-right), which also guarantees a\n# later pair's \"old\" text is never a string an earlier pair just produced\n# (verified against the source diff).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"60\u00f74=15, 0\", \"56\u00f74=14, 0\"),\n  @(\"72\u00f79=8, 0\", \"36\u00f75=7, 1\"),\n  @(\"64\u00f72=32, 0\", \"59\u00f76=9, 5\"),\n  @(\"73\u00f76=12, 1\", \"24\u00f76=4, 0\"),\n  @(\"27\u00f77=3, 6\", \"28\u00f73=9, 1\"),\n  @(\"93\u00f76=15, 3\", \"40\u00f78=5, 0\"),\n  @(\"51\u00f75=10, 1\", \"72\u00f74=18, 0\"),\n  @(\"81\u00f75=16, 1\", \"30\u00f79=3, 3\"),\n  @(\"72\u00f75=14, 2\", \"46\u00f78=5, 6\"),\n  @(\"83\u00f74=20, 3\", \"10\u00f72=5, 0\"),\n  @(\"58\u00f72=29, 0\", \"41\u00f75=8, 1\"),\n  @(\"89\u00f72=44, 1\", \"39\u00f72=19, 1\"),\n  @(\"47\u00f78=5, 7\", \"50\u00f74=12, 2\"),\n  @(\"94\u00f72=47, 0\", \"49\u00f77=7, 0\"),\n  @(\"26\u00f75=5, 1\", \"56\u00f77=8, 0\"),\n  @(\"97\u00f72=48, 1\", \"83\u00f77=11, 6\"),\n  @(\"29\u00f72=14, 1\", \"46\u00f76=7, 4\"),\n  @(\"32\u00f78=4, 0\", \"51\u00f76=8, 3\"),\n  @(\"48\u00f75=9, 3\", \"85\u00f76=14, 1\"),\n  @(\"54\u00f72=27, 0\", \"81\u00f77=11, 4\"),\n  @(\"53\u00f76=8, 5\", \"58\u00f78=7, 2\"),\n  @(\"86\u00f72=43, 0\", \"34\u00f74=8, 2\"),\n  @(\"38\u00f73=12, 2\", \"93\u00f76=15, 3\"),\n  @(\"87\u00f78=10, 7\", \"29\u00f78=3, 5\"),\n  @(\"36\u00f76=6, 0\", \"85\u00f74=21, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
